# Auto-generated edit script: replicate commit
# "إضافة حدث جديد في Card24 by F at 2025-12-28 19:15:56"
#
# The underlying data pipeline re-exports the whole workbook from its
# backing store whenever an event is added. For this commit that means:
#  1) Card24: previously-blank cells in rows 2-23 are (re)written with the
#     literal placeholder text "nan" (pandas NaN -> str), a brand-new row 24
#     is appended for the new event (card=24, Event="C"), and the sheet
#     dimension grows to A1:P24.
#  2) DATA BASE: cells A2:E11 that had been carrying the same "nan"
#     placeholder are now written out blank.

$wb = $excel.ActiveWorkbook

# ---- Card24 sheet ----
$card24 = $wb.Worksheets.Item("Card24")

# Keep all of these as plain text (the source file stores every value,
# including numbers like "24", as text/inlineStr).
$card24.Range("A2:P24").NumberFormat = "@"

$nanCells = @("E2","F2","G2","H2","I2","J2","K2","L2","N2","P2","D3","E3","F3","G3","H3","I3","J3","K3","L3","M3","N3","O3","P3","D4","E4","F4","G4","H4","I4","J4","K4","L4","M4","N4","O4","P4","D5","E5","F5","G5","H5","I5","J5","K5","L5","M5","N5","O5","P5","E6","J6","K6","M6","N6","O6","P6","E7","H7","I7","J7","K7","M7","N7","O7","P7","E8","F8","G8","H8","K8","P8","D9","E9","F9","G9","H9","I9","J9","K9","L9","M9","N9","O9","P9","D10","E10","F10","G10","H10","I10","J10","K10","L10","M10","N10","O10","P10","D11","E11","F11","G11","H11","I11","J11","K11","L11","M11","N11","O11","P11","D12","E12","F12","G12","H12","I12","J12","K12","L12","M12","N12","O12","P12","B13","C13","D13","E13","F13","G13","H13","I13","J13","K13","P13","B14","C14","D14","E14","F14","G14","H14","I14","J14","K14","M14","P14","B15","C15","D15","E15","F15","G15","H15","I15","J15","K15","M15","P15","B16","C16","D16","E16","F16","G16","H16","I16","J16","K16","P16","B17","C17","D17","E17","F17","G17","H17","I17","J17","K17","P17","B18","C18","D18","E18","F18","G18","H18","I18","J18","K18","P18","B19","C19","D19","E19","F19","G19","H19","I19","J19","K19","P19","B20","C20","D20","E20","F20","G20","H20","I20","J20","K20","M20","P20","B21","C21","D21","E21","F21","G21","H21","I21","J21","K21","M21","P21","B22","C22","D22","E22","F22","G22","H22","I22","J22","K22","M22","P22","B23","C23","D23","E23","F23","G23","H23","I23","J23","K23","P23")
foreach ($addr in $nanCells) {
    $card24.Range($addr).Value = "nan"
}

# New event row
$card24.Range("A24").Value = "24"
$card24.Range("M24").Value = "C"

# ---- DATA BASE sheet ----
$dataBase = $wb.Worksheets.Item("DATA BASE")

$clearCells = @("A2","B2","C2","D2","E2","A3","B3","C3","D3","E3","A4","B4","C4","D4","E4","A5","B5","C5","D5","E5","A6","B6","C6","D6","E6","A7","B7","C7","D7","E7","A8","B8","C8","D8","E8","A9","B9","C9","D9","E9","A10","B10","C10","D10","E10","A11","B11","C11","D11","E11")
foreach ($addr in $clearCells) {
    $dataBase.Range($addr).Value = ""
}

